$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the B-column dates up by one row (B2..B6 get the values that
# previously lived in B3..B7), then drop the now-duplicate last row.
$ws.Range("B2").Value = "Martes 04/06/2024"
$ws.Range("B3").Value = "Lunes 10/06/2024"
$ws.Range("B4").Value = "Martes 11/06/2024"
$ws.Range("B5").Value = "Lunes 24/06/2024"
$ws.Range("B6").Value = "Martes 25/06/2024"

# Remove the old row 7 entirely (was A7="2024", B7="Martes 25/06/2024").
$ws.Rows.Item(7).Delete()
